$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 data - note column D (score) is intentionally left blank to create a gap
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Lim"
$ws.Cells.Item(8, 3).Value = "Black"
$ws.Cells.Item(8, 5).Value = "male"

$f8 = $ws.Cells.Item(8, 6)
$f8.Value = 32511
$f8.NumberFormat = "mm/dd/yy;@"

$ws.Range("F8").Select()
